$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.033.46'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '1.781.03'
$ws.Range("E3").Value = '  -3.34%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = '''224.04'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '''32.40'
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("E9").Value = '  -3.77%  '
$ws.Range("D10").Value = '''0.0702'
$ws.Range("E10").Value = '  -2.16%  '
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '2.042.37'
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").Value = '1.780.04'
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("E14").Value = '  -4.17%  '
$ws.Range("D15").Value = '34.011.31'
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("D16").Value = '''0.620'
$ws.Range("E16").Value = '  -5.04%  '
$ws.Range("D17").Value = '''4.13'
$ws.Range("E17").Value = '  -5.07%  '
$ws.Range("D18").Value = '''67.58'
$ws.Range("E18").Value = '  -3.25%  '
$ws.Range("D19").Value = '''242.06'
$ws.Range("E19").Value = '  -4.73%  '
$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  -3.58%  '
$ws.Range("D21").Value = '''1.00'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '''10.66'
$ws.Range("E22").Value = '  -6.10%  '
$ws.Range("E23").Value = '  -5.49%  '
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("D25").Value = '''159.57'
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = '''16.24'
$ws.Range("E26").Value = '  -3.93%  '
$ws.Range("D27").Value = '''7.01'
$ws.Range("E27").Value = '  -3.43%  '
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.0513'
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.21'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '''3.65'
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("D33").Value = '''3.49'
$ws.Range("E33").Value = '  -4.13%  '
$ws.Range("D34").Value = '''1.80'
$ws.Range("E34").Value = '  -7.72%  '
$ws.Range("D35").Value = '1.391.42'
$ws.Range("E35").Value = '  -3.77%  '
$ws.Range("D36").Value = '''0.641'
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").Value = '''1.04'
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("E38").Value = '  -4.67%  '
$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").Value = '''2.35'
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''2.20'
$ws.Range("E40").Value = '  +1.31%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.69'
$ws.Range("E41").Value = '  -4.34%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''78.26'
$ws.Range("E42").Value = '  -5.98%  '
$ws.Range("D43").Value = '''0.906'
$ws.Range("E43").Value = '  -7.66%  '
$ws.Range("D44").Value = '0.0₆0142'
$ws.Range("E44").Value = '  +11.00%  '
$ws.Range("D45").Value = '''1.07'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("D47").Value = '''106.97'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '''5.85'
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("D49").Value = '1.939.40'
$ws.Range("E49").Value = '  -3.02%  '
$ws.Range("D50").Value = '''12.17'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("E51").Value = '  -0.21%  '
